# Update the depth/temp slope comparison table with the corrected
# statistics from the meeting with Breckie.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 2: "Deep - Mid"
$t.Cell(2, 2).Range.Text = "0.004"
$t.Cell(2, 3).Range.Text = "0.002"
$t.Cell(2, 4).Range.Text = "22"
$t.Cell(2, 5).Range.Text = "2.136"
$t.Cell(2, 6).Range.Text = "0.105"

# Row 3: "Deep - Shallow"
$t.Cell(3, 3).Range.Text = "0.002"
$t.Cell(3, 4).Range.Text = "22"
$t.Cell(3, 5).Range.Text = "1.254"
$t.Cell(3, 6).Range.Text = "0.435"

# Row 4: "Mid - Shallow"
$t.Cell(4, 3).Range.Text = "0.002"
$t.Cell(4, 4).Range.Text = "22"
$t.Cell(4, 5).Range.Text = "-0.825"
$t.Cell(4, 6).Range.Text = "0.692"
